$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Architect. Design Phase Defects")
$ws.Activate()

$ws.Range("I3").Value = "Pirlea Silvia-Cristina"
$ws.Range("J3").Value = 235

$ws.Range("C10").Value = "A01"
$ws.Range("C11").Value = "A02"
$ws.Range("C12").Value = "A03"
$ws.Range("C13").Value = "A04"
$ws.Range("C14").Value = "A05"
$ws.Range("C15").Value = "A06"
$ws.Range("C16").Value = "A07"
$ws.Range("C17").Value = "A08"
$ws.Range("C18").Value = "A09"

$ws.Range("E10").Value = "Organizararea programului este clara, exista pachete separate pentru fiecare strat al arhitecturii"
$ws.Range("E11").Value = "Partitionarea si layering-ul este corect"
$ws.Range("E12").Value = "Arhitectura permite realizarea tuturor cerintelor"

$ws.Range("D13").Value = "PizzaService"
$ws.Range("E13").Value = "La nivelul proiectului exista un singur service responabil de business layer si acesta incorporeaza toate subsistemele"

$ws.Range("E14").Value = "Nu exista la nivelul aplicatiei o strategie de gestiune a erorilor, nu exista clase separate pentru un anumit tip de eroare. Sunt tratate doar erorile IO"

$ws.Range("E15").Value = "MVC model este incorporat in proiect"

$ws.Range("D16").Value = "PizzaService"
$ws.Range("E16").Value = "Se acupa atat cu plata cat si cu afisarea meniului. Din numele clasei nu iti dai seama care este scopul serviciului, o denumire mai buna ar fi PizzaManagementOrdersService"

$ws.Range("E17").Value = "Exista descriere la clase ?? Nu cred"

$ws.Range("E18").Value = "Relatiile 1 to many nu sunt bine scrise intre PaymentRepository si Payment, intre MenuGUIController si OrdersGUI ar trebui sa fie 1:10, "

$ws.Range("C19").Value = "A10"
$ws.Range("E19").Value = "The key entity classes are consistent with business and model layers"

$ws.Range("I10").Select()
